# Add the newly-opened "divisional round" betting lines to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @(20, "SF_SEA", 46.5, 6.5),
    @(20, "LA_CHI", 51.5, -4.5),
    @(20, "BUF_DEN", 46.5, -1.5)
)

$startRow = 264
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Application.ActiveWindow.ScrollRow = 249
$ws.Range("E264").Select()
